$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pairs = @(
    ,@('69+2=71', '99-12=87')
    ,@('21+6=27', '17-10=7')
    ,@('70-20=50', '68-16=52')
    ,@('2+85=87', '35+7=42')
    ,@('94+0=94', '27+31=58')
    ,@('2+26=28', '66-41=25')
    ,@('2+37=39', '94-60=34')
    ,@('20+11=31', '61-4=57')
    ,@('91-66=25', '64+20=84')
    ,@('87-32=55', '45+16=61')
    ,@('58+23=81', '8-4=4')
    ,@('81+4=85', '40+5=45')
    ,@('55-47=8', '68-56=12')
    ,@('81-42=39', '66-11=55')
    ,@('92-51=41', '50-2=48')
    ,@('22+19=41', '9+75=84')
    ,@('51-50=1', '73+5=78')
    ,@('17+0=17', '14+45=59')
    ,@('16+63=79', '68-46=22')
    ,@('2+94=96', '82-30=52')
    ,@('15+42=57', '34-15=19')
    ,@('40+45=85', '79-31=48')
    ,@('38+3=41', '9+21=30')
    ,@('27+66=93', '6+41=47')
    ,@('13+36=49', '49-22=27')
    ,@('18+25=43', '56+5=61')
    ,@('59-53=6', '55-34=21')
    ,@('99-49=50', '86-39=47')
    ,@('7+42=49', '41+7=48')
    ,@('17+27=44', '62-30=32')
    ,@('30+65=95', '77-24=53')
    ,@('54+16=70', '7+81=88')
    ,@('19+80=99', '44+14=58')
    ,@('60-51=9', '45-43=2')
    ,@('30-8=22', '92+2=94')
    ,@('94-66=28', '90-24=66')
    ,@('61-38=23', '39+37=76')
    ,@('23+35=58', '43-22=21')
    ,@('83-16=67', '1+20=21')
    ,@('41+41=82', '79-36=43')
    ,@('75-71=4', '65+20=85')
    ,@('13+36=49', '28+53=81')
    ,@('31+19=50', '36+6=42')
    ,@('20+45=65', '79-56=23')
    ,@('48-18=30', '32+3=35')
    ,@('47+15=62', '15+14=29')
    ,@('59-15=44', '36-15=21')
    ,@('64+5=69', '21+57=78')
    ,@('41+32=73', '21-20=1')
    ,@('14+30=44', '21-6=15')
    ,@('15+83=98', '57+8=65')
    ,@('16+9=25', '13+48=61')
    ,@('60-33=27', '97-77=20')
    ,@('63-56=7', '59+8=67')
    ,@('61-49=12', '92-32=60')
    ,@('87-64=23', '27+11=38')
    ,@('18+47=65', '18-1=17')
    ,@('25+43=68', '43+3=46')
    ,@('29+11=40', '97-0=97')
    ,@('72-18=54', '82-61=21')
    ,@('80-12=68', '59-36=23')
    ,@('72+11=83', '29+14=43')
    ,@('21+31=52', '12-5=7')
    ,@('69+4=73', '31+35=66')
    ,@('33-29=4', '69+15=84')
    ,@('7+38=45', '49-14=35')
    ,@('30+56=86', '13+6=19')
    ,@('7+68=75', '35+41=76')
    ,@('9+90=99', '31+20=51')
    ,@('74-32=42', '13+41=54')
    ,@('31+55=86', '23-0=23')
    ,@('74-16=58', '33+16=49')
    ,@('93-87=6', '65-59=6')
    ,@('89-28=61', '12+22=34')
    ,@('28+35=63', '65+6=71')
    ,@('88-66=22', '20-13=7')
    ,@('85-42=43', '90+7=97')
    ,@('41+40=81', '61+5=66')
    ,@('77-33=44', '0+15=15')
    ,@('9-5=4', '1+51=52')
    ,@('39-26=13', '99-50=49')
    ,@('57-41=16', '20+48=68')
    ,@('51+25=76', '97-97=0')
    ,@('48+44=92', '40+21=61')
    ,@('43+45=88', '37+50=87')
    ,@('56+39=95', '40+25=65')
    ,@('50-39=11', '41+58=99')
    ,@('70+1=71', '61-6=55')
    ,@('22+11=33', '92-12=80')
    ,@('53+42=95', '70+19=89')
    ,@('23+2=25', '40+35=75')
    ,@('25+69=94', '77-38=39')
    ,@('9+85=94', '33+38=71')
    ,@('97+1=98', '64-35=29')
    ,@('76-67=9', '57+29=86')
    ,@('34+42=76', '84-81=3')
    ,@('3+51=54', '99-4=95')
    ,@('69+23=92', '35+25=60')
    ,@('19-10=9', '36-33=3')
    ,@('58+26=84', '44+53=97')
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $pair = $pairs[$idx]
        $old = $pair[0]
        $new = $pair[1]
        $cell = $t.Cell($r, $c)
        $cellRange = $cell.Range
        $ok = $cellRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1)
        if (-not $ok) {
            Write-Host "MISS row=$r col=$c old=$old new=$new"
        }
        $idx = $idx + 1
    }
}
Write-Host "done idx=$idx"
